$wb = $excel.ActiveWorkbook

# --- Sheet "MSCdtRPbQL": update the "Change in Perc Share" label ---
$ws2 = $wb.Worksheets.Item("MSCdtRPbQL")
$ws2.Range("A2").Value = "Change in Perc Share (dimensionless)"
$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 30
$ws2.Range("A2").Select()

# --- Sheet "About": insert an explanatory sentence before the methodology note ---
$ws1 = $wb.Worksheets.Item("About")
$ws1.Rows.Item(10).Insert()
$ws1.Range("A10").Value = "This variable measures how a rebate program influences market shares of rebate-qualifying and non-qualifying components."
$ws1.Range("A10").Font.Bold = $false
$ws1.Range("A10").Select()
